# Updated cryptos list - Price and Volume(1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.978.41"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "'1.764.83"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'322.05"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'0.9952"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.4239"
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("D8").Value = "'0.3592"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "'44.21"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").Value = "'0.07462"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "'1.104"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.9966"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'21.51"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'6.098"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "'7.320"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "'1.787.75"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "'91.52"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "'0.06370"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "'0.9952"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").Value = "'5.973"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "'28.022.10"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "'11.29"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'2.150"
$ws.Range("E25").Value = "  -6.48%  "
$ws.Range("D26").Value = "'159.47"
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "'1.986.21"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "'2.163"
$ws.Range("E29").Value = "  -5.68%  "
$ws.Range("D30").Value = "'125.42"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").Value = "'1.174"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'5.674"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "'0.09013"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").Value = "'3.498"
$ws.Range("E34").Value = "  -4.02%  "
$ws.Range("D35").Value = "'12.65"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.02316"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'5.056"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'0.06070"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "'0.2100"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").Value = "'0.6393"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "'1.190"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "'0.9947"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D45").Value = "'13.52"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'0.5957"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'3.693"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").Value = "'1.989"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "'123.44"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").Value = "'1.167"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").Value = "'0.06875"
$ws.Range("E51").Value = "  +0.03%  "

# Rows 43 and 44 swap position (WEMIXTOKEN now ranks above FraxShare)
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.401"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'7.868"
$ws.Range("E44").Value = "  -0.84%  "
